$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data to append: date serial 43933 (2020-04-12) for the four health regions
$newRows = @(
    @{ Date = 43933; Region = "Helse Midt-Norge"; Admissions = 5 },
    @{ Date = 43933; Region = "Helse Nord";        Admissions = 5 },
    @{ Date = 43933; Region = "Helse Sør-Øst";     Admissions = 40 },
    @{ Date = 43933; Region = "Helse Vest";        Admissions = 9 }
)

$startRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = $row.Date
    $dateCell.NumberFormat = "yyyy-mm-dd"

    $ws.Cells.Item($r, 2).Value = $row.Region
    $ws.Cells.Item($r, 3).Value = $row.Admissions
}
